# Reshape the little 3x5 table into a 5x4 table: the header row becomes a real
# header ("Unnamed: 0"/"2019"/"Unnamed: 1"/"2018"/"Unnamed: 2") and the two
# data columns (formerly col A & C, rows 2 & 4) move to columns B & D across
# three consecutive rows (2, 3, 4). Column A's border/bold/center-top style
# (cellXf index 1 in the original file) is reused for every header cell.

$xlPasteFormats = -4122
$xlPasteValues  = -4163

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the cells that don't survive in their old spot ---
$ws.Range("A2").Clear()
$ws.Range("C2").Clear()
$ws.Range("A3").Clear()
$ws.Range("A4").Clear()
$ws.Range("C4").Clear()
$ws.Range("A5").Clear()

# --- Header row: stamp A1's existing format onto B1:E1 first ---
$ws.Range("A1").Copy()
$ws.Range("B1:E1").PasteSpecial($xlPasteFormats)

# Headers that aren't number-like text can just be assigned directly.
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("C1").Value = "Unnamed: 1"
$ws.Range("E1").Value = "Unnamed: 2"

# "2019"/"2018" look like numbers, so a plain .Value assignment would store
# them as numeric 2019/2018 instead of text. Write each one with a leading
# quote-prefix on an out-of-the-way scratch cell (forces text storage),
# paste just the resulting text *value* into place, then re-stamp the
# header format on top (a values-only paste doesn't disturb the type).
$ws.Range("Z100").Value = "'2019"
$ws.Range("Z100").Copy()
$ws.Range("B1").PasteSpecial($xlPasteValues)

$ws.Range("Z100").Value = "'2018"
$ws.Range("Z100").Copy()
$ws.Range("D1").PasteSpecial($xlPasteValues)

$ws.Range("Z100").Clear()

$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial($xlPasteFormats)
$ws.Range("D1").PasteSpecial($xlPasteFormats)

# --- Data rows ---
$ws.Range("B2").Value = 45.4
$ws.Range("D2").Value = 45.2

$ws.Range("B3").Value = 45.3
$ws.Range("D3").Value = 45.1

$ws.Range("B4").Value = 47
$ws.Range("D4").Value = 46.8
